$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.174.41'
$ws.Range("E2").Value = '  -3.28%  '
$ws.Range("D3").Value = '3.139.97'
$ws.Range("E3").Value = '  -2.34%  '
$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = '1.00'
$r.Style = "Normal"
$ws.Range("E4").Value = '  +0.09%  '
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '605.97'
$r.Style = "Normal"
$ws.Range("E5").Value = '  -0.25%  '
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '147.63'
$r.Style = "Normal"
$ws.Range("E6").Value = '  -5.50%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").Value = '3.139.24'
$ws.Range("E9").Value = '  -3.53%  '
$ws.Range("E10").Value = '  -5.21%  '
$ws.Range("E11").Value = '  -2.47%  '
$ws.Range("E12").Value = '  -5.21%  '
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = '0.0000258'
$r.Style = "Normal"
$ws.Range("E13").Value = '  -3.71%  '
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '36.57'
$r.Style = "Normal"
$ws.Range("E14").Value = '  -4.69%  '
$ws.Range("D15").Value = '3.657.49'
$ws.Range("E15").Value = '  -2.29%  '
$ws.Range("D16").Value = '64.217.48'
$ws.Range("E16").Value = '  -3.33%  '
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = '0.113'
$r.Style = "Normal"
$ws.Range("E17").Value = '  -0.01%  '
$ws.Range("D18").Value = '3.065.80'
$ws.Range("E18").Value = '  -4.66%  '
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = '6.96'
$r.Style = "Normal"
$ws.Range("E19").Value = '  -4.21%  '
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '480.42'
$r.Style = "Normal"
$ws.Range("E20").Value = '  -5.24%  '
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '14.56'
$r.Style = "Normal"
$ws.Range("E21").Value = '  -4.68%  '
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = '0.708'
$r.Style = "Normal"
$ws.Range("E22").Value = '  -2.85%  '
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '7.72'
$r.Style = "Normal"
$ws.Range("E23").Value = '  -4.04%  '
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = '13.76'
$r.Style = "Normal"
$ws.Range("E24").Value = '  -5.66%  '
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = '83.77'
$r.Style = "Normal"
$ws.Range("E25").Value = '  -1.46%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("E27").Value = '  -2.42%  '
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = '8.52'
$r.Style = "Normal"
$ws.Range("E28").Value = '  -5.06%  '
$ws.Range("E29").Value = '  -4.95%  '
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = '0.123'
$r.Style = "Normal"
$ws.Range("E30").Value = '  -11.77%  '
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '6.85'
$r.Style = "Normal"
$ws.Range("E31").Value = '  -1.16%  '
$ws.Range("E32").Value = '  -5.66%  '
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = '1.00'
$r.Style = "Normal"
$ws.Range("E33").Value = '  -0.08%  '
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = '26.66'
$r.Style = "Normal"
$ws.Range("E34").Value = '  -5.64%  '
$ws.Range("E35").Value = '  -4.96%  '
$ws.Range("E36").Value = '  -5.11%  '
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = '54.49'
$r.Style = "Normal"
$ws.Range("E37").Value = '  -1.71%  '
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '3.13'
$r.Style = "Normal"
$ws.Range("E38").Value = '  +2.81%  '
$ws.Range("D39").Value = '0.0₃0736'
$ws.Range("E39").Value = '  -4.76%  '
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = '452.22'
$r.Style = "Normal"
$ws.Range("E40").Value = '  -9.81%  '
$ws.Range("E41").Value = '  -4.65%  '
$ws.Range("E42").Value = '  -5.82%  '
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '8.43'
$r.Style = "Normal"
$ws.Range("E43").Value = '  -3.46%  '
$ws.Range("D44").Value = '2.869.09'
$ws.Range("E44").Value = '  -1.91%  '
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '0.270'
$r.Style = "Normal"
$ws.Range("E45").Value = '  -8.41%  '
$ws.Range("E46").Value = '  -6.82%  '
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = '26.60'
$r.Style = "Normal"
$ws.Range("E47").Value = '  -5.36%  '
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = '2.31'
$r.Style = "Normal"
$ws.Range("E49").Value = '  -3.37%  '
$ws.Range("E50").Value = '  -2.49%  '
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = '120.14'
$r.Style = "Normal"
$ws.Range("E51").Value = '  -1.29%  '

